# Fruta / hortaliza, semanal
# Insert a new weekly record above the existing row 711, pushing all
# subsequent rows down by one (old row 711 -> 712, ..., old row 807 -> 808).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 711 - shifts rows 711:807 down to 712:808
# and grows the sheet dimension to A1:R808.
$ws.Rows.Item(711).Insert()

# Populate the newly inserted row 711 with the new observation.
$ws.Range("A711").Value = 3
$ws.Range("B711").Value = "Femacal de La Calera"
$ws.Range("C711").Value = "Coquimbo"
$ws.Range("D711").Value = 45131
$ws.Range("E711").Value = 5
$ws.Range("F711").Value = 100112032
$ws.Range("G711").Value = "Zapallo italiano"
$ws.Range("H711").Value = "Sin especificar"
$ws.Range("I711").Value = "Primera"
$ws.Range("J711").Value = 80
$ws.Range("K711").Value = 13500
$ws.Range("L711").Value = 14000
$ws.Range("M711").Value = 13750
$ws.Range("N711").Value = "`$/caja 60 unidades"
$ws.Range("O711").Value = "Región de Arica y Parinacota"
$ws.Range("P711").Value = 229
$ws.Range("Q711").Value = 60
$ws.Range("R711").Value = "Hortaliza"

# Match the date-number-format style used by the rest of column D.
$ws.Range("D711").NumberFormat = $ws.Range("D710").NumberFormat
